# Automated map update - adds 5 new PEBCOM case rows (95-99) to the sheet,
# matching the "Actualización automática del mapa" data refresh.
#
# Columns: A Caso | B F.De Reclamo | C Direccion | D Comuna | E OT |
#          F Proveedor Asignado | G Estado | H Observaciones |
#          I Attachments | J Tipo de tarea | K Equipo | L Tipo de Elemento |
#          M Coordenada_X | N Coordenada_Y | O Operacion | P Zona | Q PD | R N2
#
# Columns A, B, D, E contain values that look numeric/date-like (case ids,
# comuna numbers, OT numbers, dd/m/yyyy dates) but must be stored as literal
# text, exactly as they already are for every pre-existing row in the sheet.
# Pre-formatting each of those cells as Text ("@") before assigning the
# value is what makes Excel keep the literal text instead of silently
# re-interpreting it as a number or a date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 95 : Caso 7451 - CORVALAN 3698 ----
$ws.Range("A95").NumberFormat = "@"
$ws.Range("A95").Value = '7451'
$ws.Range("B95").NumberFormat = "@"
$ws.Range("B95").Value = '10/8/2025'
$ws.Range("C95").Value = 'CORVALAN 3698'
$ws.Range("D95").NumberFormat = "@"
$ws.Range("D95").Value = '8'
$ws.Range("E95").NumberFormat = "@"
$ws.Range("E95").Value = '810259135'
$ws.Range("F95").Value = 'PEBCOM'
$ws.Range("G95").Value = 'Pendiente'
$ws.Range("H95").Value = 'Cambiar'
$ws.Range("I95").Value = 1
$ws.Range("J95").Value = 'Cambio'
$ws.Range("K95").Value = 'Sin equipos'
$ws.Range("L95").Value = 'Terminal'
$ws.Range("M95").Value = -58.46777
$ws.Range("N95").Value = -34.671445
$ws.Range("O95").Value = 'Boedo'
$ws.Range("P95").Value = 'Capital Sur'
$ws.Range("Q95").Value = 'PAV-?'
$ws.Range("R95").Value = 'Fuera de Poligono OVL'

# ---- Row 96 : Caso 7453 - GUARDIA NACIONAL 2616 ----
$ws.Range("A96").NumberFormat = "@"
$ws.Range("A96").Value = '7453'
$ws.Range("B96").NumberFormat = "@"
$ws.Range("B96").Value = '10/8/2025'
$ws.Range("C96").Value = 'GUARDIA NACIONAL 2616'
$ws.Range("D96").NumberFormat = "@"
$ws.Range("D96").Value = '8'
$ws.Range("E96").NumberFormat = "@"
$ws.Range("E96").Value = '810259140'
$ws.Range("F96").Value = 'PEBCOM'
$ws.Range("G96").Value = 'Pendiente'
$ws.Range("H96").Value = 'Picada'
$ws.Range("I96").Value = 1
$ws.Range("J96").Value = 'Cambio'
$ws.Range("K96").Value = 'Sin equipos'
$ws.Range("L96").Value = 'Pasante'
$ws.Range("M96").Value = -58.476253
$ws.Range("N96").Value = -34.660961
$ws.Range("O96").Value = 'Boedo'
$ws.Range("P96").Value = 'Capital Sur'
$ws.Range("Q96").Value = 'PAV-P'
$ws.Range("R96").Value = 'Fuera de Poligono OVL'

# ---- Row 97 : Caso 7477 - GAONA AV. 5130 ----
$ws.Range("A97").NumberFormat = "@"
$ws.Range("A97").Value = '7477'
$ws.Range("B97").NumberFormat = "@"
$ws.Range("B97").Value = '10/8/2025'
$ws.Range("C97").Value = 'GAONA AV. 5130'
$ws.Range("D97").NumberFormat = "@"
$ws.Range("D97").Value = '10'
$ws.Range("E97").NumberFormat = "@"
$ws.Range("E97").Value = '810259143'
$ws.Range("F97").Value = 'PEBCOM'
$ws.Range("G97").Value = 'Pendiente'
$ws.Range("H97").Value = 'Picada'
$ws.Range("I97").Value = 1
$ws.Range("J97").Value = 'Cambio'
$ws.Range("K97").Value = 'Sin equipos'
$ws.Range("L97").Value = 'Pasante'
$ws.Range("M97").Value = -58.493913
$ws.Range("N97").Value = -34.62931
$ws.Range("O97").Value = 'Devoto'
$ws.Range("P97").Value = 'Capital Norte'
$ws.Range("Q97").Value = 'DEV-M'
$ws.Range("R97").Value = 'ARATO-25058.PO.2DEV'

# ---- Row 98 : Caso 7481 - BACACAY 2455 ----
$ws.Range("A98").NumberFormat = "@"
$ws.Range("A98").Value = '7481'
$ws.Range("B98").NumberFormat = "@"
$ws.Range("B98").Value = '10/8/2025'
$ws.Range("C98").Value = 'BACACAY 2455'
$ws.Range("D98").NumberFormat = "@"
$ws.Range("D98").Value = '7'
$ws.Range("E98").NumberFormat = "@"
$ws.Range("E98").Value = '810259148'
$ws.Range("F98").Value = 'PEBCOM'
$ws.Range("G98").Value = 'Pendiente'
$ws.Range("H98").Value = 'Picada'
$ws.Range("I98").Value = 1
$ws.Range("J98").Value = 'Cambio'
$ws.Range("K98").Value = 'Sin equipos'
$ws.Range("L98").Value = 'Pasante'
$ws.Range("M98").Value = -58.464662
$ws.Range("N98").Value = -34.626638
$ws.Range("O98").Value = 'Boedo'
$ws.Range("P98").Value = 'Capital Sur'
$ws.Range("Q98").Value = 'NRA-E'
$ws.Range("R98").Value = 'Fuera de Poligono OVL'

# ---- Row 99 : Caso -634 - Curapaligue 1127 ----
$ws.Range("A99").NumberFormat = "@"
$ws.Range("A99").Value = '-634'
$ws.Range("B99").NumberFormat = "@"
$ws.Range("B99").Value = '10/8/2025'
$ws.Range("C99").Value = 'Curapaligue 1127'
$ws.Range("D99").NumberFormat = "@"
$ws.Range("D99").Value = '7'
$ws.Range("E99").NumberFormat = "@"
$ws.Range("E99").Value = 'Pendiente ADM'
$ws.Range("F99").Value = 'PEBCOM'
$ws.Range("G99").Value = 'Pendiente'
$ws.Range("H99").Value = 'Colocar columna donde se marca en la foto pasante 150 o 200'
$ws.Range("I99").Value = 1
$ws.Range("J99").Value = 'Cambio'
$ws.Range("K99").Value = 'Sin equipos'
$ws.Range("L99").Value = 'Pasante'
$ws.Range("M99").Value = -58.446624
$ws.Range("N99").Value = -34.635851
$ws.Range("O99").Value = 'Boedo'
$ws.Range("P99").Value = 'Capital Sur'
$ws.Range("Q99").Value = 'PPT-M'
$ws.Range("R99").Value = 'Fuera de Poligono OVL'
